$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.631.63"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "1.587.65"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.65"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  -2.64%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.246"
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.52"
$ws.Range("E10").Value = "  -4.27%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "1.810.08"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "1.585.07"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  -4.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.67"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "26.614.43"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "0.0₃0725"
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.81"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.85"
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.56"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.25"
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0506"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.656"
$ws.Range("E33").Value = "  +19.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.89"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").Value = "1.303.76"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.791"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.49"
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("D45").Value = "1.723.59"
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.59"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.838"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0977"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.49"
$ws.Range("E51").Value = "  -1.66%  "
